# Update cryptos list with freshly scraped price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.614.20"
$ws.Range("E2").Value = "  +5.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.067.62"
$ws.Range("E3").Value = "  +6.10%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.30"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.94"
$ws.Range("E6").Value = "  +5.06%  "
$ws.Range("E7").Value = "  +21.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.060.40"
$ws.Range("E8").Value = "  +6.26%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +10.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.177"
$ws.Range("E11").Value = "  +5.53%  "
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.81"
$ws.Range("E13").Value = "  +18.08%  "
$ws.Range("E14").Value = "  +10.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.707.26"
$ws.Range("E15").Value = "  +5.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.085.05"
$ws.Range("E16").Value = "  +6.85%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.39"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.44"
$ws.Range("E18").Value = "  +5.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.24"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.606.56"
$ws.Range("E21").Value = "  +5.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.73"
$ws.Range("E22").Value = "  +7.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.18"
$ws.Range("E23").Value = "  +21.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.61"
$ws.Range("E24").Value = "  +7.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.13"
$ws.Range("E25").Value = "  +8.63%  "
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.08"
$ws.Range("E28").Value = "  +6.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.13"
$ws.Range("E29").Value = "  +6.95%  "
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.29"
$ws.Range("E31").Value = "  +17.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.71"
$ws.Range("E32").Value = "  +5.67%  "
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "678.92"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "68.02"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.61"
$ws.Range("E36").Value = "  +13.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.24"
$ws.Range("E37").Value = "  +7.70%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0865"
$ws.Range("E38").Value = "  +3.48%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.432"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +4.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  +9.29%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  +5.69%  "
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.160"
$ws.Range("E45").Value = "  +15.58%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.22"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.88"
$ws.Range("E47").Value = "  +17.79%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +5.87%  "
$ws.Range("E51").Value = "  +4.12%  "
